# Update "想去人数" (interest count) values for a few events whose numbers
# changed between data pulls (regenerated gh-pages data at 456a3b4).
# The same events appear both on the "展览" (exhibitions) sheet and the
# "全部类型" (all types) sheet, so both must be kept in sync.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 665
$ws1.Range("F9").Value = 7948
$ws1.Range("F13").Value = 383

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 665
$ws4.Range("F11").Value = 7948
$ws4.Range("F17").Value = 383
